$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("HomePage", "SignOutBtn", "xpath", "//a[@class='logout hidden-sm-down']"),
    @("HomePage", "SearchCtlgTxtField", "name", "s"),
    @("HomePage", "MyStoreLogo", "xpath", "//img[@class='logo img-fluid']"),
    @("HomePage", "ClothesMenuLink", "xpath", "//li[@id='category-3']"),
    @("HomePage", "AccessoriesMenuLink", "xpath", "//li[@id='category-6']"),
    @("HomePage", "ArtMenuLink", "xpath", "//li[@id='category-9']"),
    @("HomePage", "SubMenuLink", "xpath", "//a[@class='dropdown-item dropdown-submenu']")
)

$row = 23
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $row++
}

$ws.Columns("D:D").AutoFit() | Out-Null

$ws.Range("D33").Select()
